$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("08-09-2021", 1000000, 1665000, 1500000, 900000, 600000, 1.49),
    @("09-09-2021", 1200000, 2402000, 1200000, 593000, 607000, 1.47),
    @("10-09-2021", 1200000, 2190000, 1200000, 610000, 590000, 1.45),
    @("13-09-2021", 500000, 1465000, 750000, 530000, 220000, 1.4),
    @("14-09-2021", 600000, 965000, 900000, 745000, 155000, 1.49),
    @("15-09-2021", 2200000, 2875000, 2200000, 1445000, 755000, 1.48),
    @("16-09-2021", 2600000, 2805000, 2600000, 1575000, 1025000, 1.49)
)

$startRow = 174
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]

    # Column A holds a day-month-year text label (e.g. "08-09-2021"). Excel's
    # automatic type inference would otherwise turn some of these into real
    # dates, so force a text number format before assigning the value, then
    # restore the Normal style so no extra formatting is left on the cell.
    $cellA = $ws.Cells.Item($row, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $vals[0]
    $cellA.Style = "Normal"

    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
    $ws.Cells.Item($row, 6).Value = $vals[5]
    $ws.Cells.Item($row, 7).Value = $vals[6]
}
